# Automatische test-sync: 2025-06-26 23:01:50
# Appends a new Logs row (test mail #1 reply about opening hours) and bumps
# the matching Dashboard category count.

$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append row 30 -------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A30").Value = "Wanneer zijn jullie open?"
$wsLogs.Range("B30").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C30").Value = "Testmail #1: Wanneer zijn jullie open?"
$wsLogs.Range("D30").Value = "Openingstijden / Locatie"
$wsLogs.Range("E30").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$wsLogs.Range("F30").Value = "2025-06-26 23:01:23"
$wsLogs.Range("G30").Value = "Ja"
$wsLogs.Range("H30").Value = "Nee"
$wsLogs.Range("I30").Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I) from row 29 to row 30,
# same as Excel does automatically when a table-like range grows.
$oldRanges = @("D2:D29", "G2:G29", "H2:H29", "I2:I29")
$newRanges = @("D2:D30", "G2:G30", "H2:H30", "I2:I30")

for ($j = 0; $j -lt $oldRanges.Count; $j++) {
    $fcs = $wsLogs.Range($oldRanges[$j]).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($wsLogs.Range($newRanges[$j]))
    }
}

# --- "Dashboard" sheet: bump the "Openingstijden / Locatie" tally ----------
$wsDashboard = $wb.Worksheets.Item("Dashboard")
$wsDashboard.Range("B3").Value = 7

Write-Output "Logs row 30 added; Dashboard B3 updated to 7."
